$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bulk Upload")

$ws.Range("A2").Value = "Proctor70909"
$ws.Range("B2").Value = "Automation70909"
$ws.Range("C2").Value = "proctorautomation70909@gmail.com"
$ws.Range("D2").Value = "'70909"
$ws.Range("D2").Style = "Normal"

$ws.Range("A3").Value = "Proctor05280"
$ws.Range("B3").Value = "Automation05280"
$ws.Range("C3").Value = "proctorautomation05280@gmail.com"
$ws.Range("D3").Value = "'05280"
$ws.Range("D3").Style = "Normal"

$ws.Range("A4").Value = "Proctor75082"
$ws.Range("B4").Value = "Automation75082"
$ws.Range("C4").Value = "proctorautomation75082@gmail.com"
$ws.Range("D4").Value = "'75082"
$ws.Range("D4").Style = "Normal"

$ws.Range("A5").Value = "Proctor15299"
$ws.Range("B5").Value = "Automation15299"
$ws.Range("C5").Value = "proctorautomation15299@gmail.com"
$ws.Range("D5").Value = "'15299"
$ws.Range("D5").Style = "Normal"

$ws.Range("A6").Value = "Proctor11174"
$ws.Range("B6").Value = "Automation11174"
$ws.Range("C6").Value = "proctorautomation11174@gmail.com"
$ws.Range("D6").Value = "'11174"
$ws.Range("D6").Style = "Normal"

$ws.Range("A7").Value = "Proctor07661"
$ws.Range("B7").Value = "Automation07661"
$ws.Range("C7").Value = "proctorautomation07661@gmail.com"
$ws.Range("D7").Value = "'07661"
$ws.Range("D7").Style = "Normal"

$ws.Range("A8").Value = "Proctor57545"
$ws.Range("B8").Value = "Automation57545"
$ws.Range("C8").Value = "proctorautomation57545@gmail.com"
$ws.Range("D8").Value = "'57545"
$ws.Range("D8").Style = "Normal"

$ws.Range("A9").Value = "Proctor87841"
$ws.Range("B9").Value = "Automation87841"
$ws.Range("C9").Value = "proctorautomation87841@gmail.com"
$ws.Range("D9").Value = "'87841"
$ws.Range("D9").Style = "Normal"

$ws.Range("A10").Value = "Proctor37284"
$ws.Range("B10").Value = "Automation37284"
$ws.Range("C10").Value = "proctorautomation37284@gmail.com"
$ws.Range("D10").Value = "'37284"
$ws.Range("D10").Style = "Normal"

$ws.Range("A11").Value = "Proctor28592"
$ws.Range("B11").Value = "Automation28592"
$ws.Range("C11").Value = "proctorautomation28592@gmail.com"
$ws.Range("D11").Value = "'28592"
$ws.Range("D11").Style = "Normal"

$ws.Range("A12").Value = "Proctor96408"
$ws.Range("B12").Value = "Automation96408"
$ws.Range("C12").Value = "proctorautomation96408@gmail.com"
$ws.Range("D12").Value = "'96408"
$ws.Range("D12").Style = "Normal"

$ws.Range("A13").Value = "Proctor42276"
$ws.Range("B13").Value = "Automation42276"
$ws.Range("C13").Value = "proctorautomation42276@gmail.com"
$ws.Range("D13").Value = "'42276"
$ws.Range("D13").Style = "Normal"

$ws.Range("A14").Value = "Proctor64899"
$ws.Range("B14").Value = "Automation64899"
$ws.Range("C14").Value = "proctorautomation64899@gmail.com"
$ws.Range("D14").Value = "'64899"
$ws.Range("D14").Style = "Normal"

$ws.Range("A15").Value = "Proctor21299"
$ws.Range("B15").Value = "Automation21299"
$ws.Range("C15").Value = "proctorautomation21299@gmail.com"
$ws.Range("D15").Value = "'21299"
$ws.Range("D15").Style = "Normal"

$ws.Range("A16").Value = "Proctor66618"
$ws.Range("B16").Value = "Automation66618"
$ws.Range("C16").Value = "proctorautomation66618@gmail.com"
$ws.Range("D16").Value = "'66618"
$ws.Range("D16").Style = "Normal"

$ws.Range("A17").Value = "Proctor84933"
$ws.Range("B17").Value = "Automation84933"
$ws.Range("C17").Value = "proctorautomation84933@gmail.com"
$ws.Range("D17").Value = "'84933"
$ws.Range("D17").Style = "Normal"

$ws.Range("A18").Value = "Proctor52628"
$ws.Range("B18").Value = "Automation52628"
$ws.Range("C18").Value = "proctorautomation52628@gmail.com"
$ws.Range("D18").Value = "'52628"
$ws.Range("D18").Style = "Normal"

$ws.Range("A19").Value = "Proctor11072"
$ws.Range("B19").Value = "Automation11072"
$ws.Range("C19").Value = "proctorautomation11072@gmail.com"
$ws.Range("D19").Value = "'11072"
$ws.Range("D19").Style = "Normal"

$ws.Range("A20").Value = "Proctor59923"
$ws.Range("B20").Value = "Automation59923"
$ws.Range("C20").Value = "proctorautomation59923@gmail.com"
$ws.Range("D20").Value = "'59923"
$ws.Range("D20").Style = "Normal"

$ws.Range("A21").Value = "Proctor91490"
$ws.Range("B21").Value = "Automation91490"
$ws.Range("C21").Value = "proctorautomation91490@gmail.com"
$ws.Range("D21").Value = "'91490"
$ws.Range("D21").Style = "Normal"

$ws.Range("A22").Value = "Proctor32281"
$ws.Range("B22").Value = "Automation32281"
$ws.Range("C22").Value = "proctorautomation32281@gmail.com"
$ws.Range("D22").Value = "'32281"
$ws.Range("D22").Style = "Normal"

$ws.Range("A23").Value = "Proctor14105"
$ws.Range("B23").Value = "Automation14105"
$ws.Range("C23").Value = "proctorautomation14105@gmail.com"
$ws.Range("D23").Value = "'14105"
$ws.Range("D23").Style = "Normal"

$ws.Range("A24").Value = "Proctor86192"
$ws.Range("B24").Value = "Automation86192"
$ws.Range("C24").Value = "proctorautomation86192@gmail.com"
$ws.Range("D24").Value = "'86192"
$ws.Range("D24").Style = "Normal"

$ws.Range("A25").Value = "Proctor73249"
$ws.Range("B25").Value = "Automation73249"
$ws.Range("C25").Value = "proctorautomation73249@gmail.com"
$ws.Range("D25").Value = "'73249"
$ws.Range("D25").Style = "Normal"

$ws.Range("A26").Value = "Proctor13507"
$ws.Range("B26").Value = "Automation13507"
$ws.Range("C26").Value = "proctorautomation13507@gmail.com"
$ws.Range("D26").Value = "'13507"
$ws.Range("D26").Style = "Normal"

$ws.Range("A27").Value = "Proctor20531"
$ws.Range("B27").Value = "Automation20531"
$ws.Range("C27").Value = "proctorautomation20531@gmail.com"
$ws.Range("D27").Value = "'20531"
$ws.Range("D27").Style = "Normal"

$ws.Range("A28").Value = "Proctor82597"
$ws.Range("B28").Value = "Automation82597"
$ws.Range("C28").Value = "proctorautomation82597@gmail.com"
$ws.Range("D28").Value = "'82597"
$ws.Range("D28").Style = "Normal"

$ws.Range("A29").Value = "Proctor75542"
$ws.Range("B29").Value = "Automation75542"
$ws.Range("C29").Value = "proctorautomation75542@gmail.com"
$ws.Range("D29").Value = "'75542"
$ws.Range("D29").Style = "Normal"

$ws.Range("A30").Value = "Proctor58492"
$ws.Range("B30").Value = "Automation58492"
$ws.Range("C30").Value = "proctorautomation58492@gmail.com"
$ws.Range("D30").Value = "'58492"
$ws.Range("D30").Style = "Normal"

$ws.Range("A31").Value = "Proctor26968"
$ws.Range("B31").Value = "Automation26968"
$ws.Range("C31").Value = "proctorautomation26968@gmail.com"
$ws.Range("D31").Value = "'26968"
$ws.Range("D31").Style = "Normal"
